{"js": "// Cover letter typo fixes (\"Main Point\" paragraph):\n//   1. \"the vision modular robots\"               -> \"the vision of modular robots\"\n//   2. \"meet the needs of a perceived,\"           -> \"meet the needs of perceived,\"\n//   3. \"unknown environment. This results marks\"  -> \"unknown environments. This results marks\"\n\n// 1. Insert \"of \" between \"the vision \" and \"modular robots reconfiguring\".\n{\n  const results = context.document.body.search(\n    \"the vision modular robots reconfiguring\",\n    { matchCase: true, matchWholeWord: false }\n  );\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"the vision of modular robots reconfiguring\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// 2. Remove the \"a \" before \"perceived,\" (second occurrence, after \"meet the needs of\").\n{\n  const results = context.document.body.search(\n    \"meet the needs of a perceived,\",\n    { matchCase: true, matchWholeWord: false }\n  );\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"meet the needs of perceived,\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// 3. Add the missing \"s\" to the second \"unknown environment\" (the one right before\n//    \"This results marks a milestone...\").\n{\n  const results = context.document.body.search(\n    \"unknown environment. This results marks\",\n    { matchCase: true, matchWholeWord: false }\n  );\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"unknown environments. This results marks\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n", "ps1": "# Cover letter typo fixes (\"Main Point\" paragraph):\n#   1. \"the vision modular robots\"               -> \"the vision of modular robots\"\n#   2. \"meet the needs of a perceived,\"           -> \"meet the needs of perceived,\"\n#   3. \"unknown environment. This results marks\"  -> \"unknown environments. This results marks\"\n\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceOne = 2 (Replace:=2 substitutes only the single found match)\n\n# 1. Insert \"of \" between \"the vision \" and \"modular robots reconfiguring\".\n$find1 = $d.Content.Find\n$find1.Execute(\n    \"the vision modular robots reconfiguring\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"the vision of modular robots reconfiguring\",\n    2\n) | Out-Null\n\n# 2. Remove the \"a \" before \"perceived,\" (the one following \"meet the needs of\").\n$find2 = $d.Content.Find\n$find2.Execute(\n    \"meet the needs of a perceived,\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"meet the needs of perceived,\",\n    2\n) | Out-Null\n\n# 3. Add the missing \"s\" to the second \"unknown environment\" (the one right before\n#    \"This results marks a milestone...\").\n$find3 = $d.Content.Find\n$find3.Execute(\n    \"unknown environment. This results marks\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"unknown environments. This results marks\",\n    2\n) | Out-Null\n"}
